$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (C) column date value from 45172 to 45175 for every
# data row (rows 2 through 420).
$ws.Range("C2:C420").Value = 45175
